$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 "detail task" relabeling: "프로세스 설계" -> "기능 구상", "프로세스 구상" -> "기능 구상/기술 결정"
$ws.Range("B4").Value = "기능 구상"
$ws.Range("C4").Value = "기능 구상/기술 결정"

# Duration swap between row6 (퍼블리싱) and row7 (출시 및 유지보수)
$ws.Range("E6").Value = 29
$ws.Range("E7").Value = 21

# Align F7's formatting (date number format + vertical-center alignment) with the rest of column F
$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# View state: zoom + scroll position + active selection
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("I7").Select() | Out-Null
